$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append the new "// Ngày 12/12/2020" run right after "SW7A3R".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("SW7A3R", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" // Ngày 12/12/2020")

$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Color = 16777215
$r.Font.Size = 13.5
$r.Style = "Strong"

# ---------------------------------------------------------------------------
# 2. Drop a collapsed "_GoBack" bookmark right at the very end of the body
#    (after the text just inserted). A temporary trailing character is used
#    so the bookmark does not land exactly on the paragraph-end boundary
#    (where it would otherwise get reseated to the paragraph start).
# ---------------------------------------------------------------------------
$r.Collapse(0)
$endPos = $r.End

$r.InsertAfter("Z")
$bmRange = $d.Range($endPos, $endPos + 1)
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($endPos, $endPos + 1)
$delRange.Delete()
